# The document contains four "<id>p048r_aN</id>" tags (each built from
# three separate runs: "<id>", the bare id text, and "</id>"). The commit
# renames each identifier from "p048r_aN" to "p048r_N" (dropping the "a"),
# which collapses the three runs into a single run holding the full
# "<id>p048r_N</id>" text.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p048r_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p048r_1</id>", 2)
$d.Content.Find.Execute("<id>p048r_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p048r_2</id>", 2)
$d.Content.Find.Execute("<id>p048r_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p048r_3</id>", 2)
$d.Content.Find.Execute("<id>p048r_a4</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p048r_4</id>", 2)
